$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before D (old D shifts to F, E to G, etc.)
$ws.Columns("D:E").Insert()

# Apply number formats to the newly inserted columns, matching column F
$ws.Range("D8:E102").NumberFormat = $ws.Range("F8").NumberFormat
$ws.Range("D7:E7").NumberFormat = $ws.Range("F7").NumberFormat
$ws.Range("D38:E38").NumberFormat = $ws.Range("F38").NumberFormat
$ws.Range("D80:E80").NumberFormat = $ws.Range("F80").NumberFormat

# Populate new columns D and E with the added Q4-2018 / Q3-2018 data (and some blanks stay blank)
$ws.Cells.Item(7,4).Value = 43465
$ws.Cells.Item(7,5).Value = 43373
$ws.Cells.Item(8,4).Value = 289900
$ws.Cells.Item(8,5).Value = 207200
$ws.Cells.Item(9,4).Value = 55400
$ws.Cells.Item(9,5).Value = 50000
$ws.Cells.Item(10,4).Value = 234500
$ws.Cells.Item(10,5).Value = 157200
$ws.Cells.Item(12,4).Value = "NA"
$ws.Cells.Item(12,5).Value = "NA"
$ws.Cells.Item(13,4).Value = 0
$ws.Cells.Item(13,5).Value = 0
$ws.Cells.Item(14,4).Value = 0
$ws.Cells.Item(14,5).Value = 31200
$ws.Cells.Item(15,4).Value = 72500
$ws.Cells.Item(15,5).Value = 70500
$ws.Cells.Item(17,4).Value = 141800
$ws.Cells.Item(17,5).Value = 170600
$ws.Cells.Item(18,4).Value = 148100
$ws.Cells.Item(18,5).Value = 36600
$ws.Cells.Item(20,4).Value = 2800
$ws.Cells.Item(20,5).Value = -1200
$ws.Cells.Item(21,4).Value = 223400
$ws.Cells.Item(21,5).Value = 105900
$ws.Cells.Item(22,4).Value = 14500
$ws.Cells.Item(22,5).Value = 10300
$ws.Cells.Item(23,4).Value = 136400
$ws.Cells.Item(23,5).Value = 25100
$ws.Cells.Item(24,4).Value = -7700
$ws.Cells.Item(24,5).Value = "NA"
$ws.Cells.Item(25,4).Value = 0
$ws.Cells.Item(25,5).Value = 0
$ws.Cells.Item(26,4).Value = 144100
$ws.Cells.Item(26,5).Value = 25100
$ws.Cells.Item(27,4).Value = 136700
$ws.Cells.Item(27,5).Value = 17800
$ws.Cells.Item(28,4).Value = 0
$ws.Cells.Item(28,5).Value = 0
$ws.Cells.Item(29,4).Value = 0
$ws.Cells.Item(29,5).Value = "NA"
$ws.Cells.Item(30,4).Value = 0
$ws.Cells.Item(30,5).Value = 0
$ws.Cells.Item(31,4).Value = 0
$ws.Cells.Item(31,5).Value = 0
$ws.Cells.Item(32,4).Value = -2800
$ws.Cells.Item(32,5).Value = 1200
$ws.Cells.Item(33,4).Value = 136700
$ws.Cells.Item(33,5).Value = 17800
$ws.Cells.Item(34,4).Value = 0
$ws.Cells.Item(34,5).Value = 0
$ws.Cells.Item(35,4).Value = 136700
$ws.Cells.Item(35,5).Value = 17800
$ws.Cells.Item(38,4).Value = 43465
$ws.Cells.Item(38,5).Value = 43373
$ws.Cells.Item(41,4).Value = 64500
$ws.Cells.Item(41,5).Value = 45900
$ws.Cells.Item(42,4).Value = 0
$ws.Cells.Item(42,5).Value = 0
$ws.Cells.Item(43,4).Value = 146200
$ws.Cells.Item(43,5).Value = 146900
$ws.Cells.Item(44,4).Value = 17600
$ws.Cells.Item(44,5).Value = 18800
$ws.Cells.Item(45,4).Value = 77400
$ws.Cells.Item(45,5).Value = 13900
$ws.Cells.Item(46,4).Value = 305700
$ws.Cells.Item(46,5).Value = 225400
$ws.Cells.Item(47,4).Value = 0
$ws.Cells.Item(47,5).Value = 0
$ws.Cells.Item(48,4).Value = 3122900
$ws.Cells.Item(48,5).Value = 2922000
$ws.Cells.Item(49,4).Value = 0
$ws.Cells.Item(49,5).Value = 0
$ws.Cells.Item(50,4).Value = 0
$ws.Cells.Item(50,5).Value = 0
$ws.Cells.Item(51,4).Value = 0
$ws.Cells.Item(51,5).Value = 0
$ws.Cells.Item(52,4).Value = 27000
$ws.Cells.Item(52,5).Value = 6800
$ws.Cells.Item(53,4).Value = 0
$ws.Cells.Item(53,5).Value = 0
$ws.Cells.Item(54,4).Value = 3455500
$ws.Cells.Item(54,5).Value = 3154200
$ws.Cells.Item(57,4).Value = 77900
$ws.Cells.Item(57,5).Value = 44800
$ws.Cells.Item(58,4).Value = 0
$ws.Cells.Item(58,5).Value = 0
$ws.Cells.Item(59,4).Value = 252100
$ws.Cells.Item(59,5).Value = 281900
$ws.Cells.Item(60,4).Value = 330000
$ws.Cells.Item(60,5).Value = 326800
$ws.Cells.Item(61,4).Value = 1297800
$ws.Cells.Item(61,5).Value = 1065100
$ws.Cells.Item(62,4).Value = 48000
$ws.Cells.Item(62,5).Value = 39900
$ws.Cells.Item(63,4).Value = 0
$ws.Cells.Item(63,5).Value = 0
$ws.Cells.Item(64,4).Value = 0
$ws.Cells.Item(64,5).Value = 0
$ws.Cells.Item(65,4).Value = 0
$ws.Cells.Item(65,5).Value = 0
$ws.Cells.Item(66,4).Value = 1766600
$ws.Cells.Item(66,5).Value = 1607000
$ws.Cells.Item(68,4).Value = 0
$ws.Cells.Item(68,5).Value = 0
$ws.Cells.Item(69,4).Value = 0
$ws.Cells.Item(69,5).Value = 0
$ws.Cells.Item(70,4).Value = 0
$ws.Cells.Item(70,5).Value = 0
$ws.Cells.Item(71,4).Value = 0
$ws.Cells.Item(71,5).Value = 0
$ws.Cells.Item(72,4).Value = -236300
$ws.Cells.Item(72,5).Value = -373000
$ws.Cells.Item(73,4).Value = 0
$ws.Cells.Item(73,5).Value = 0
$ws.Cells.Item(74,4).Value = 0
$ws.Cells.Item(74,5).Value = 0
$ws.Cells.Item(75,4).Value = 0
$ws.Cells.Item(75,5).Value = 0
$ws.Cells.Item(76,4).Value = 1688900
$ws.Cells.Item(76,5).Value = 1547200
$ws.Cells.Item(77,4).Value = 0
$ws.Cells.Item(77,5).Value = 0
$ws.Cells.Item(80,4).Value = 43465
$ws.Cells.Item(80,5).Value = 43373
$ws.Cells.Item(81,4).Value = 136700
$ws.Cells.Item(81,5).Value = 17800
$ws.Cells.Item(83,4).Value = 72500
$ws.Cells.Item(83,5).Value = 70500
$ws.Cells.Item(84,4).Value = 0
$ws.Cells.Item(84,5).Value = 0
$ws.Cells.Item(85,4).Value = 0
$ws.Cells.Item(85,5).Value = 0
$ws.Cells.Item(86,4).Value = 0
$ws.Cells.Item(86,5).Value = 0
$ws.Cells.Item(87,4).Value = 0
$ws.Cells.Item(87,5).Value = 0
$ws.Cells.Item(88,4).Value = 0
$ws.Cells.Item(88,5).Value = 0
$ws.Cells.Item(89,4).Value = 189200
$ws.Cells.Item(89,5).Value = 165100
$ws.Cells.Item(91,4).Value = -43500
$ws.Cells.Item(91,5).Value = -42700
$ws.Cells.Item(92,4).Value = 0
$ws.Cells.Item(92,5).Value = 0
$ws.Cells.Item(93,4).Value = 0
$ws.Cells.Item(93,5).Value = 0
$ws.Cells.Item(94,4).Value = -294700
$ws.Cells.Item(94,5).Value = -727000
$ws.Cells.Item(96,4).Value = 0
$ws.Cells.Item(96,5).Value = 0
$ws.Cells.Item(97,4).Value = 0
$ws.Cells.Item(97,5).Value = 0
$ws.Cells.Item(98,4).Value = 0
$ws.Cells.Item(98,5).Value = 0
$ws.Cells.Item(99,4).Value = 0
$ws.Cells.Item(99,5).Value = 0
$ws.Cells.Item(100,4).Value = 136500
$ws.Cells.Item(100,5).Value = 471400
$ws.Cells.Item(101,4).Value = 0
$ws.Cells.Item(101,5).Value = 0
$ws.Cells.Item(102,4).Value = 31000
$ws.Cells.Item(102,5).Value = -90500

# A handful of historical quarters were restated in this edit; correct those shifted cells
$ws.Cells.Item(91,6).Value = -42200
$ws.Cells.Item(91,7).Value = -37300
$ws.Cells.Item(91,8).Value = -40300
$ws.Cells.Item(91,9).Value = -38800
$ws.Cells.Item(91,10).Value = -20900
$ws.Cells.Item(94,8).Value = -227200
$ws.Cells.Item(94,9).Value = -227200
$ws.Cells.Item(102,8).Value = 66800
$ws.Cells.Item(102,9).Value = -115600

Write-Output "edit complete"